$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 234
$ws.Range("F6").Value = 79
$ws.Range("F7").Value = 804
$ws.Range("F8").Value = 456
$ws.Range("F9").Value = 65
$ws.Range("F10").Value = 278
$ws.Range("F11").Value = 101
$ws.Range("F12").Value = 214
$ws.Range("F13").Value = 15
$ws.Range("F15").Value = 6409
$ws.Range("F17").Value = 63
$ws.Range("F19").Value = 7374
$ws.Range("F21").Value = 31
$ws.Range("F22").Value = 3340
$ws.Range("F23").Value = 770
$ws.Range("F24").Value = 848
$ws.Range("F25").Value = 4497
$ws.Range("F26").Value = 339
$ws.Range("F27").Value = 177
$ws.Range("F29").Value = 1401
$ws.Range("F30").Value = 135
$ws.Range("F31").Value = 46
$ws.Range("F33").Value = 1086
$ws.Range("F34").Value = 1511
$ws.Range("F35").Value = 2106

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 39
$ws.Range("F5").Value = 72

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1194
$ws.Range("F4").Value = 66

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1194
$ws.Range("F5").Value = 66
$ws.Range("F8").Value = 234
$ws.Range("F9").Value = 79
$ws.Range("F10").Value = 804
$ws.Range("F11").Value = 456
$ws.Range("F12").Value = 65
$ws.Range("F13").Value = 278
$ws.Range("F15").Value = 101
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 15
$ws.Range("F19").Value = 6409
$ws.Range("F21").Value = 63
$ws.Range("F23").Value = 7374
$ws.Range("F25").Value = 31
$ws.Range("F26").Value = 3340
$ws.Range("F27").Value = 770
$ws.Range("F28").Value = 848
$ws.Range("F29").Value = 4497
$ws.Range("F30").Value = 339
$ws.Range("F31").Value = 39
$ws.Range("F32").Value = 177
$ws.Range("F34").Value = 1401
$ws.Range("F35").Value = 135
$ws.Range("F36").Value = 46
$ws.Range("F38").Value = 1086
$ws.Range("F39").Value = 1511
$ws.Range("F41").Value = 2106
$ws.Range("F43").Value = 72
